$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save off row 2's current values for columns D, L, M, N, O, P, S
$row2_D = $ws.Range("D2").Value()
$row2_L = $ws.Range("L2").Value()
$row2_M = $ws.Range("M2").Value()
$row2_N = $ws.Range("N2").Value()
$row2_O = $ws.Range("O2").Value()
$row2_P = $ws.Range("P2").Value()
$row2_S = $ws.Range("S2").Value()

# Save off row 5's current values for columns D, L, M, N, O, P, S
$row5_D = $ws.Range("D5").Value()
$row5_L = $ws.Range("L5").Value()
$row5_M = $ws.Range("M5").Value()
$row5_N = $ws.Range("N5").Value()
$row5_O = $ws.Range("O5").Value()
$row5_P = $ws.Range("P5").Value()
$row5_S = $ws.Range("S5").Value()

# Write row 5's old values into row 2
$ws.Range("D2").Value = $row5_D
$ws.Range("L2").Value = $row5_L
$ws.Range("M2").Value = $row5_M
$ws.Range("N2").Value = $row5_N
$ws.Range("O2").Value = $row5_O
$ws.Range("P2").Value = $row5_P
$ws.Range("S2").Value = $row5_S

# Write row 2's old values into row 5
$ws.Range("D5").Value = $row2_D
$ws.Range("L5").Value = $row2_L
$ws.Range("M5").Value = $row2_M
$ws.Range("N5").Value = $row2_N
$ws.Range("O5").Value = $row2_O
$ws.Range("P5").Value = $row2_P
$ws.Range("S5").Value = $row2_S
